$d = $word.ActiveDocument

$find = "For my research project, I extended a research project in Unity to leverage and explore the VR space in 3D graph exploration."
$repl = "For my research project, I extended the functionality of an existing Unity application to leverage and explore the VR space in 3D graph exploration."

$d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $repl, 2)

try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
} catch {
}
